$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.183.16"
$ws.Range("E2").Value = "  +1.70%  "
$ws.Range("D3").Value = "3.214.26"
$ws.Range("E3").Value = "  +1.05%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "'604.92"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +4.31%  "
$ws.Range("D6").Value = "'153.96"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.27%  "
$ws.Range("D7").Value = "'1.00"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("D8").Value = "3.212.18"
$ws.Range("E8").Value = "  +0.94%  "
$ws.Range("D9").Value = "'0.535"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.28%  "
$ws.Range("E10").Value = "  -1.04%  "
$ws.Range("E11").Value = "  -0.97%  "
$ws.Range("E12").Value = "  +1.02%  "
$ws.Range("D13").Value = "'0.0000273"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.15%  "
$ws.Range("D14").Value = "'38.64"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.81%  "
$ws.Range("D15").Value = "3.739.91"
$ws.Range("E15").Value = "  +1.02%  "
$ws.Range("B16").Value = "WrappedBTC"
$ws.Range("C16").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D16").Value = "66.222.02"
$ws.Range("E16").Value = "  +1.67%  "
$ws.Range("B17").Value = "Polkadot"
$ws.Range("C17").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D17").Value = "'7.47"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +3.86%  "
$ws.Range("D18").Value = "3.213.48"
$ws.Range("E18").Value = "  +1.14%  "
$ws.Range("E19").Value = "  -0.32%  "
$ws.Range("D20").Value = "'512.32"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.47%  "
$ws.Range("D21").Value = "'15.56"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +4.29%  "
$ws.Range("E22").Value = "  +0.21%  "
$ws.Range("B23").Value = "InternetComputer(DFINITY)"
$ws.Range("C23").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D23").Value = "'15.28"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.65%  "
$ws.Range("B24").Value = "Uniswap"
$ws.Range("C24").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D24").Value = "'8.04"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +2.88%  "
$ws.Range("D25").Value = "'85.23"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.28%  "
$ws.Range("E26").Value = "  -0.12%  "
$ws.Range("E27").Value = "  +2.63%  "
$ws.Range("D28").Value = "'9.21"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +2.26%  "
$ws.Range("D29").Value = "'2.26"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +2.97%  "
$ws.Range("D30").Value = "'2.87"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +2.91%  "
$ws.Range("D31").Value = "'6.83"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +8.00%  "
$ws.Range("D32").Value = "'28.17"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.70%  "
$ws.Range("D33").Value = "'1.22"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.98%  "
$ws.Range("E34").Value = "  +0.10%  "
$ws.Range("E35").Value = "  +0.28%  "
$ws.Range("D36").Value = "'55.37"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.73%  "
$ws.Range("D37").Value = "'0.0908"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.61%  "
$ws.Range("D38").Value = "'481.83"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.21%  "
$ws.Range("E39").Value = "  -0.30%  "
$ws.Range("E40").Value = "  -5.47%  "
$ws.Range("E41").Value = "  +2.04%  "
$ws.Range("D42").Value = "'0.298"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +3.66%  "
$ws.Range("E43").Value = "  -0.54%  "
$ws.Range("D44").Value = "2.945.48"
$ws.Range("E44").Value = "  -4.15%  "
$ws.Range("E45").Value = "  +1.88%  "
$ws.Range("D46").Value = "0.0₃0642"
$ws.Range("E46").Value = "  +4.52%  "
$ws.Range("D47").Value = "'28.87"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.07%  "
$ws.Range("E49").Value = "  -0.22%  "
$ws.Range("E50").Value = "  +2.25%  "
$ws.Range("B51").Value = "Arweave"
$ws.Range("C51").Value = "https://coinranking.com/coin/7XWg41D1+arweave-ar"
$ws.Range("D51").Value = "'33.88"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +4.56%  "
